$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The uploaded images were converted from .jpg to .png, so the filename
# list in column A needs its extensions updated to match (row 1 "photo"
# header is untouched).
$ws.Range("A2").Value = "Images/photoI.png"
$ws.Range("A3").Value = "Images/photoII.png"
$ws.Range("A4").Value = "Images/photoIII.png"
$ws.Range("A5").Value = "Images/photoIV.png"
$ws.Range("A6").Value = "Images/photoV.png"
$ws.Range("A7").Value = "Images/photoVI.png"
$ws.Range("A8").Value = "Images/photoVII.png"

# Re-font the sheet from Tahoma to Calibri.
$ws.Cells.Font.Name = "Calibri"

# Leave the cursor on the last touched cell, like the author did before
# saving/uploading.
$ws.Range("D7").Select()
